# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns on the per-language sheets (zh-cn, de-de) now that a handback has happened, and
# flips the "Status" shown on the Overview + language sheets from "Ready for handoff" to
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$hyperlinkUrlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e9251d2fd2344319e672556d3cbdbbad0eeeea6/e2e/a.md"

$hyperlinkColor = 15570276   # RGB(100,149,237) == FF6495ED, matching the existing custom "HyperLink" cell style

# ---------------------------------------------------------------------------
# Overview sheet: widen the zh-cn / de-de status columns and update the status
# text shown for both languages (both rows).
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.2
$overview.Columns.Item(6).ColumnWidth = 29.2

# ---------------------------------------------------------------------------
# Helper: apply the handback update to a language sheet (zh-cn / de-de).
# ---------------------------------------------------------------------------
function Update-LanguageSheet {
    param($ws, $handbackFile, $handbackDateTime)

    # Status column (C) gets wider + new text on both data rows.
    $ws.Columns.Item(3).ColumnWidth = 29.2
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Handback File column (J) gets wider to fit the long file name.
    $ws.Columns.Item(10).ColumnWidth = 39.2

    # Latest Target File (I): now links to a.md, matching the handoff file link.
    $ws.Range("I2").Value = "a.md"
    $ws.Range("I3").Value = "a.md"
    $ws.Hyperlinks.Add($ws.Range("I2"), $hyperlinkUrlA, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $hyperlinkUrlA, "", "", "a.md")
    foreach ($ref in @("I2", "I3")) {
        $cell = $ws.Range($ref)
        $cell.Style = "HyperLink"
        $cell.Font.Underline = $true
        $cell.Font.Color = $hyperlinkColor
    }

    # Latest Handback File (J): the generated handback xliff file name.
    $ws.Range("J2").Value = $handbackFile
    $ws.Range("J3").Value = $handbackFile

    # Latest Handback DateTime (K): stored as plain text, like the rest of the
    # date-ish columns in this report.
    $ws.Range("K2").Value = $handbackDateTime
    $ws.Range("K3").Value = $handbackDateTime
}

$zhcn = $wb.Worksheets.Item("zh-cn")
Update-LanguageSheet $zhcn "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-08-27 02:36:34"

$dede = $wb.Worksheets.Item("de-de")
Update-LanguageSheet $dede "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-08-27 02:36:40"
